$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 80; this shifts the existing rows 80..106
# down to 81..107 (matching the rest of the diff, which is otherwise an
# unchanged re-numbering of those rows).
$ws.Rows.Item(80).Insert()

# Populate the newly inserted row 80 with its data. Most fields mirror
# what used to be row 80 (now row 81); only Fecha (D), Precio mínimo (K),
# Precio máximo (L), Precio promedio ponderado (M), Región (O) and
# Volumen kilos (P) differ from that row.
$ws.Range("A80").Value = 10
$ws.Range("B80").Value = "Vega Modelo de Temuco"
$ws.Range("C80").Value = "La Araucanía"
$ws.Range("D80").Value = 44917
$ws.Range("E80").Value = 9
$ws.Range("F80").Value = 100112022
$ws.Range("G80").Value = "Arveja Verde"
$ws.Range("H80").Value = "Sin especificar"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 65
$ws.Range("K80").Value = 25000
$ws.Range("L80").Value = 25000
$ws.Range("M80").Value = 25000
$ws.Range("N80").Value = "$/saco 25 kilos"
$ws.Range("O80").Value = "Provincia de Cautín"
$ws.Range("P80").Value = 1000
$ws.Range("Q80").Value = 25
$ws.Range("R80").Value = "Hortaliza"
